$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1.45
$ws.Range("K5").Value = 2.2
$ws.Range("L5").Value = 2.05
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("AC5").Value = 8
$ws.Range("AL5").Value = 34
$ws.Range("AP5").Value = 41
$ws.Range("AT5").Value = 2.63
$ws.Range("BC5").Value = 126

# Row 6
$ws.Range("G6").Value = 1.27
$ws.Range("H6").Value = 5.25
$ws.Range("I6").Value = 13
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 6.5
$ws.Range("AC6").Value = 11
$ws.Range("AT6").Value = 3.25

# Row 7
$ws.Range("G7").Value = 4.1
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 1.75
$ws.Range("L7").Value = 2.38
$ws.Range("U7").Value = 1.73
$ws.Range("V7").Value = 2
$ws.Range("AG7").Value = 8
$ws.Range("AH7").Value = 9
$ws.Range("AP7").Value = 29
$ws.Range("AU7").Value = 8
$ws.Range("AX7").Value = 9

# Row 12
$ws.Range("G12").Value = 1.73
$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 2.38
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("AC12").Value = 8
$ws.Range("AO12").Value = 9
$ws.Range("AR12").Value = 51
$ws.Range("AU12").Value = 9

# Row 14
$ws.Range("G14").Value = 5.25
$ws.Range("H14").Value = 3.6
$ws.Range("I14").Value = 1.55
$ws.Range("J14").Value = 6
$ws.Range("L14").Value = 2.2
$ws.Range("N14").Value = 9.5
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.4
$ws.Range("R14").Value = 1.8
$ws.Range("W14").Value = 13
$ws.Range("X14").Value = 29
$ws.Range("Y14").Value = 19
$ws.Range("Z14").Value = 67
$ws.Range("AB14").Value = 51
$ws.Range("AD14").Value = 7.5
$ws.Range("AG14").Value = 6
$ws.Range("AH14").Value = 7
$ws.Range("AJ14").Value = 11
$ws.Range("AN14").Value = 7
$ws.Range("AO14").Value = 34
$ws.Range("AP14").Value = 41
$ws.Range("AQ14").Value = 126
$ws.Range("AR14").Value = 151
$ws.Range("AS14").Value = 351
$ws.Range("AW14").Value = 3.4
$ws.Range("AX14").Value = 8
$ws.Range("AZ14").Value = 26

# Row 15
$ws.Range("H15").Value = 2.88
$ws.Range("J15").Value = 4.75
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 2.2
$ws.Range("U15").Value = 2.25
$ws.Range("V15").Value = 1.57
$ws.Range("AE15").Value = 21
$ws.Range("AH15").Value = 8.5
$ws.Range("AR15").Value = 151
$ws.Range("AS15").Value = 451
$ws.Range("AT15").Value = 2.2
$ws.Range("BD15").Value = 126
